# Atualizando meios de saida
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "27/02/2023"
$ws.Range("B2").Value = "07:57"
$ws.Range("C2").Value = "Guilherme"
$ws.Range("D2").Value = 88359
$ws.Range("E2").Value = 88358
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "R$ 7.152.446,82"
$ws.Range("I2").Value = "R$ 28,00"
$ws.Range("J2").Value = "R$ 0,00"
$ws.Range("K2").Value = "R$ 7.152.474,82"
$ws.Range("L2").Value = "2023-02-25 23:59:41"
$ws.Range("M2").Value = "2023-02-22 00:00:00"
$ws.Range("N2").Value = "2023-02-03 11:50:27.167000"
$ws.Range("O2").Value = "2023-02-24 16:02:26.523000"
$ws.Range("P2").Value = "Não há registros a serem exibidos."
